$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph ("Existe um padrão para cMOOCs?") gains a hanging
#    indent: w:ind left=2100 leftChars=0 firstLine=700 firstLineChars=0
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 105.0
$p1.Format.FirstLineIndent = 35.0
$p1.Format.CharacterUnitLeftIndent = 0
$p1.Format.CharacterUnitFirstLineIndent = 0

# ---------------------------------------------------------------------------
# 2) Blogs paragraph: fix the botched "g" / "RSShopper" split (an artifact
#    of the old _GoBack bookmark breaking "gRSShopper" into separate runs)
#    so the line reads "...WordPress, gRSShopper, Blogger, Tumblr, Twitter"
#    as one continuous run. This also removes the old _GoBack bookmark.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Google Sites, WordPress, gRSShopper, Blogger, Tumblr, Twitter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Google Sites, WordPress, gRSShopper, Blogger, Tumblr, Twitter", 2)

# ---------------------------------------------------------------------------
# 3) Fóruns paragraph: add "Yahoo Groups, Facebook e" to the alternatives
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Google Groups, Diigo,  Fóruns em phpBB",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Google Groups, Diigo, Yahoo Groups, Facebook e  Fóruns em phpBB", 2)

# ---------------------------------------------------------------------------
# 4) Hospedagem geral de conteúdo paragraph:
#    - "Youtube (para vídeos)" -> "Youtube ou Vimeo (para vídeos)"
#    - insert "Google Docs (documentos), " before "Dropbox, ..."
#    - relocate the _GoBack bookmark to sit right before "Dropbox"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Youtube (para vídeos), Dropbox",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Youtube ou Vimeo (para vídeos), Google Docs (documentos), Dropbox", 2)

$bmRange = $d.Content
$found = $bmRange.Find.Execute(
    "Dropbox, Google Drive, OneDrive, Mega.co.nz, Amazon Cloud Drive, Box.com , Mediafire",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 5) Final content paragraph (bitnami link) switches from the generic grey
#    (#222222) to red (#FF0000) - both the run and the paragraph mark.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*bitnami*") {
        $p.Range.Font.Color = 255
    }
}
